$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Productdata")
$ws.Range("C2").Value = 0
$ws.Range("E2").Value = 1284.605555555555
$ws.Range("C3").Value = 0
$ws.Range("E3").Value = 473.4510937499999
$ws.Range("C4").Value = 0
$ws.Range("E4").Value = 110.6595486111111
$ws.Range("C5").Value = 0
$ws.Range("E5").Value = 417.2167777777777
$ws.Range("C6").Value = 0
$ws.Range("E6").Value = 482.7916909722222
$ws.Range("C7").Value = 19133
$ws.Range("E7").Value = 473.4510937499999
$ws.Range("C8").Value = 7188
$ws.Range("E8").Value = 110.6595486111111
$ws.Range("C9").Value = 22640
$ws.Range("E9").Value = 834.4335555555554
$ws.Range("C10").Value = 0
$ws.Range("E10").Value = 65.57491319444443
$ws.Range("C11").Value = 0
$ws.Range("E11").Value = 412.0635937499999
$ws.Range("C12").Value = 0
$ws.Range("E12").Value = 96.30371527777775
$ws.Range("C13").Value = 0
$ws.Range("E13").Value = 363.0034444444444
$ws.Range("C14").Value = 0
$ws.Range("E14").Value = 420.0575243055555

$ws = $wb.Worksheets.Item("ForecastedAverageDemand")
$ws.Range("G2").Value = 12144
$ws.Range("H2").Value = 1713
$ws.Range("I2").Value = 6908
$ws.Range("G3").Value = 6222
$ws.Range("H3").Value = 4756
$ws.Range("I3").Value = 9593
$ws.Range("J3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 345
$ws.Range("I4").Value = 4908
$ws.Range("J4").Value = 0
$ws.Range("G5").Value = 14665
$ws.Range("H5").Value = 7794
$ws.Range("I5").Value = 18455
$ws.Range("J5").Value = 2507
$ws.Range("G6").Value = 27062
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("G7").Value = 755
$ws.Range("H7").Value = 1729
$ws.Range("I7").Value = 4765
$ws.Range("J7").Value = 2578
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 890
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1844
$ws.Range("G9").Value = 12817
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 20427
$ws.Range("J9").Value = 3296

$ws = $wb.Worksheets.Item("ForcastedStandardDeviation")
$ws.Range("G2").Value = 303.5999999999999
$ws.Range("H2").Value = 42.82499999999999
$ws.Range("I2").Value = 172.7
$ws.Range("G3").Value = 295.5449999999999
$ws.Range("H3").Value = 225.9099999999999
$ws.Range("I3").Value = 455.6674999999998
$ws.Range("J3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 23.37374999999999
$ws.Range("I4").Value = 332.5169999999999
$ws.Range("J4").Value = 0
$ws.Range("G5").Value = 1260.823375
$ws.Range("H5").Value = 670.08915
$ws.Range("I5").Value = 1586.668625
$ws.Range("J5").Value = 215.539325
$ws.Range("G6").Value = 2770.539905
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("G7").Value = 88.44051124999999
$ws.Range("H7").Value = 202.53462775
$ws.Range("I7").Value = 558.17090875
$ws.Range("J7").Value = 301.9862755
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 116.07893975
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 240.5051290999999
$ws.Range("G9").Value = 1824.9254423575
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 2908.4615753325
$ws.Range("J9").Value = 469.2950189599999

$ws = $wb.Worksheets.Item("Capacity")
$ws.Range("B2").Value = 1156145
$ws.Range("B3").Value = 276243.75
$ws.Range("B4").Value = 64601.25
$ws.Range("B5").Value = 243960
$ws.Range("B6").Value = 282303.75
$ws.Range("B7").Value = 460406.25
$ws.Range("B8").Value = 43067.5
$ws.Range("B9").Value = 243960
$ws.Range("B10").Value = 38343.75
$ws.Range("B11").Value = 184162.5
$ws.Range("B12").Value = 21533.75
$ws.Range("B13").Value = 162640
$ws.Range("B14").Value = 376405

$ws = $wb.Worksheets.Item("ProcessingTime")
$ws.Range("B2").Value = 4
$ws.Range("C3").Value = 3
$ws.Range("D4").Value = 3
$ws.Range("E5").Value = 3
$ws.Range("F6").Value = 3
$ws.Range("G7").Value = 5
$ws.Range("H8").Value = 2
$ws.Range("I9").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M13").Value = 2
$ws.Range("N14").Value = 4
